{"js": "// Adds the reflection write-up to the (empty) document:\n//   1) \"Names: Rahul, Nick, Violet, Yacine\"\n//   2) <blank line>                              <- the original empty paragraph\n//   3) \"Each of us did one question separately and If we needed ...\"\n//   4) <blank line>\n//   5) \"From this class activity I learned about float ... ::after selector better.\"\n//\n// We build the new paragraphs (including the w:proofErr spell/grammar-check\n// markers Word leaves around \"Yacine\", \"If\" and \"the ::after\") as raw OOXML\n// and splice them in immediately before the document's first (and, in the\n// starting document, only) paragraph. Inserting at a zero-length \"Start\"\n// range with location \"Before\" keeps that original empty paragraph intact\n// (rather than merging text into it), which is what the target markup needs.\n\nconst body = context.document.body;\nconst firstParagraph = body.paragraphs.getFirst();\n\nconst newParagraphsXml =\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Names: Rahul, Nick, Violet, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Yacine</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Each of us did one question separately and </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>If</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> we needed any help we would ask help from each other. We did this for the first 5 questions but worked together on answering the last question. </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">From this class activity I learned about float and this helped me to understand </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>the ::after</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> selector better. </w:t></w:r>' +\n  '</w:p>';\n\nconst flatOpc =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphsXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nconst insertionPoint = firstParagraph.getRange(\"Start\");\ninsertionPoint.insertOoxml(flatOpc, \"Before\");\n\nawait context.sync();\n", "ps1": "# Adds the reflection write-up to the (empty) document:\n#   1) \"Names: Rahul, Nick, Violet, Yacine\"\n#   2) <blank line>                              <- the original empty paragraph\n#   3) \"Each of us did one question separately and If we needed ...\"\n#   4) <blank line>\n#   5) \"From this class activity I learned about float ... ::after selector better.\"\n#\n# The new paragraphs (including the w:proofErr spell/grammar-check markers\n# Word leaves around \"Yacine\", \"If\" and \"the ::after\") are built as raw\n# Flat-OPC OOXML and spliced in via Range.InsertXML at a zero-length range\n# collapsed to the very start of the document. That keeps the document's\n# original empty paragraph intact as paragraph #2 instead of merging text\n# into it.\n\n$d = $word.ActiveDocument\n\n$insertionPoint = $d.Paragraphs(1).Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n\n$innerBody =\n  \"<w:p>\" +\n    \"<w:r><w:t xml:space='preserve'>Names: Rahul, Nick, Violet, </w:t></w:r>\" +\n    \"<w:proofErr w:type='spellStart'/>\" +\n    \"<w:r><w:t>Yacine</w:t></w:r>\" +\n    \"<w:proofErr w:type='spellEnd'/>\" +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n    \"<w:r><w:t xml:space='preserve'>Each of us did one question separately and </w:t></w:r>\" +\n    \"<w:proofErr w:type='gramStart'/>\" +\n    \"<w:r><w:t>If</w:t></w:r>\" +\n    \"<w:proofErr w:type='gramEnd'/>\" +\n    \"<w:r><w:t xml:space='preserve'> we needed any help we would ask help from each other. We did this for the first 5 questions but worked together on answering the last question. </w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n    \"<w:r><w:t xml:space='preserve'>From this class activity I learned about float and this helped me to understand </w:t></w:r>\" +\n    \"<w:proofErr w:type='gramStart'/>\" +\n    \"<w:r><w:t>the ::after</w:t></w:r>\" +\n    \"<w:proofErr w:type='gramEnd'/>\" +\n    \"<w:r><w:t xml:space='preserve'> selector better. </w:t></w:r>\" +\n  \"</w:p>\"\n\n$flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          \"<w:body>$innerBody</w:body>\" +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$insertionPoint.InsertXML($flatOpc)\n"}
